$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: CV ID row - update SCORE, CLUSTER MUST/GOOD HAVE MATCH+SCORE ---
$ws.Range("D2").Value = "27.81"
$ws.Range("F2").Value = "sql : 1"
$ws.Range("G2").Value = "28.57"
$ws.Range("H2").Value = "requests : 3"
$ws.Range("I2").Value = "28.57"
$ws.Range("K2").Value = "0.0"

# --- Row 3: cluster match lists ---
$ws.Range("F3").Value = "analysis : 1"
$ws.Range("H3").Value = "process : 1"

# --- Row 4: cluster match lists (H4 removed, F4 replaced) ---
$ws.Range("F4").Value = "version control : 2"
$ws.Range("H4").ClearContents()

# --- Row 5: cluster match lists ---
$ws.Range("F5").Value = "databases : 2"

# --- Rows 6-8: remove stray F column entries ---
$ws.Range("F6").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("F8").ClearContents()

# --- CV KEYWORDS column (E2:E151): rewrite with updated keyword list ---
# (a handful of keywords were dropped from the analysis: coding, communication,
# windows, maintenance, code, project:9, development, training)
$keywords = @(
    'oracle dba : 1',
    'dba : 8',
    'oracle database : 28',
    'database administrator : 7',
    'administrator : 7',
    'testing : 3',
    'database : 36',
    'oracle : 33',
    'interpersonal skills : 4',
    'analytical : 1',
    'problem solving : 2',
    'ms access : 2',
    'access : 1',
    'ms sql : 2',
    'sql server : 3',
    'jdbc : 2',
    'microsoft visio : 2',
    'visio : 1',
    'sql developer : 2',
    'developer : 1',
    'toad : 2',
    'rman : 12',
    'asm : 3',
    'oem : 4',
    'grid : 6',
    'cloud : 7',
    'remedy : 1',
    'emc : 2',
    'dynatrace : 3',
    'mysql : 4',
    'encryption : 1',
    'citrix : 2',
    'solaris : 3',
    'unix : 3',
    'linux : 4',
    'database administration : 6',
    'administration : 19',
    'postgresql : 4',
    'management : 16',
    'weblogic : 2',
    'performance tuning : 9',
    'tuning : 9',
    'monitoring : 12',
    'backup : 8',
    'recovery : 9',
    'shell scripting : 3',
    'scripting : 2',
    'oracle rac : 1',
    'rac : 3',
    'installation : 6',
    'patching : 6',
    'migration : 6',
    'sql tuning : 1',
    'production : 12',
    'deployment : 1',
    'reporting : 7',
    'enterprise manager : 8',
    'manager : 7',
    'performance monitoring : 4',
    'sql : 9',
    'databases : 20',
    'communication skills : 1',
    'design : 7',
    'agile : 2',
    'customer service : 3',
    'research : 1',
    'configuration : 7',
    'legal : 1',
    'security : 3',
    'troubleshooting : 3',
    'systems : 2',
    'oracle 12c : 1',
    'hardware : 1',
    'red hat : 1',
    'operating system : 2',
    'scheduling : 5',
    'capacity planning : 2',
    'planning : 4',
    'fusion : 1',
    'middleware : 1',
    'install : 6',
    'change management : 3',
    'best practices : 1',
    'automation : 6',
    'tools : 2',
    'operations : 4',
    'oracle databases : 3',
    'documentation : 2',
    'sharepoint : 2',
    'tibco : 1',
    'shell scripts : 3',
    'tables : 5',
    'optimization : 5',
    'infrastructure : 1',
    'monitor : 1',
    'dynamic environment : 1',
    'production support : 2',
    'requests : 1',
    'business requirements : 5',
    'project management : 2',
    'analysis : 3',
    'migrations : 3',
    'server administration : 2',
    'dashboard : 1',
    'business continuity : 1',
    'process : 5',
    'version control : 2',
    'interactive : 1',
    'global : 1',
    'programmer : 1',
    'analyst : 2',
    'information system : 1',
    'stored procedures : 2',
    'debugging : 1',
    'application development : 1',
    'r&d : 3',
    'xml : 1',
    'ftp : 2',
    'business process : 3',
    'pl : 5',
    'query optimization : 6',
    'data warehouse : 3',
    'warehouse : 7',
    'contracts : 1',
    'business objects : 5',
    'basic : 3',
    'analyze : 1',
    'metadata : 1',
    'data extraction : 2',
    'sql scripts : 2',
    'unix shell : 2',
    'data warehousing : 1',
    'order management : 4',
    'data mart : 1',
    'visual basic : 2',
    'ado : 1',
    'dimensional modeling : 2',
    'modeling : 1',
    'shell : 2',
    'c : 1',
    'scheduler : 1',
    'reports : 1',
    'computer applications : 2',
    'computer science : 4',
    'foundation : 1',
    'service management : 1',
    'data protection : 1',
    'oracle e-business suite : 1',
    'ms sql server : 1',
    'pl/sql : 2'
)

for ($i = 0; $i -lt $keywords.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $keywords[$i]
}

# --- Remove now-unused trailing rows 152:160 (keyword list shrank by 9 overall) ---
$ws.Range("A152:K160").EntireRow.Delete()
